# Update LR-pair NATMI stats with refreshed TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"46.11811066666667"
$ws.Range("H2").Value = [double]"138.354332"
$ws.Range("I2").Value = [double]"0.95896098489411"
$ws.Range("J2").Value = [double]"0.9589609848941099"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.2662156666666667"
$ws.Range("N2").Value = [double]"0.7986470000000001"
$ws.Range("O2").Value = [double]"0.0009813702709097034"
$ws.Range("P2").Value = [double]"0.0009813702709097034"
$ws.Range("Q2").Value = [double]"12.27736357653378"
$ws.Range("R2").Value = [double]"110.496272188804"
$ws.Range("S2").Value = [double]"0.0009410958015373687"
$ws.Range("T2").Value = [double]"0.0009410958015373686"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"46.11811066666667"
$ws.Range("H3").Value = [double]"138.354332"
$ws.Range("I3").Value = [double]"0.95896098489411"
$ws.Range("J3").Value = [double]"0.9589609848941099"
$ws.Range("O3").Value = [double]"0.998256289001958"
$ws.Range("P3").Value = [double]"0.998256289001958"
$ws.Range("Q3").Value = [double]"12488.61491522204"
$ws.Range("R3").Value = [double]"112397.5342369984"
$ws.Range("S3").Value = [double]"0.9572888340780569"
$ws.Range("T3").Value = [double]"0.9572888340780568"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"46.11811066666667"
$ws.Range("H4").Value = [double]"138.354332"
$ws.Range("I4").Value = [double]"0.95896098489411"
$ws.Range("J4").Value = [double]"0.9589609848941099"
$ws.Range("M4").Value = [double]"0.2067996666666667"
$ws.Range("N4").Value = [double]"0.620399"
$ws.Range("O4").Value = [double]"0.000762340727132399"
$ws.Range("P4").Value = [double]"0.0007623407271323989"
$ws.Range("Q4").Value = [double]"9.537209913163112"
$ws.Range("R4").Value = [double]"85.83488921846801"
$ws.Range("S4").Value = [double]"0.0007310550145157773"
$ws.Range("T4").Value = [double]"0.0007310550145157771"
$ws.Range("I5").Value = [double]"0.002799731840346333"
$ws.Range("J5").Value = [double]"0.002799731840346333"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.2662156666666667"
$ws.Range("N5").Value = [double]"0.7986470000000001"
$ws.Range("O5").Value = [double]"0.0009813702709097034"
$ws.Range("P5").Value = [double]"0.0009813702709097034"
$ws.Range("Q5").Value = [double]"0.03584434222266666"
$ws.Range("R5").Value = [double]"0.322599080004"
$ws.Range("S5").Value = [double]"2.747573594635203E-06"
$ws.Range("T5").Value = [double]"2.747573594635203E-06"
$ws.Range("I6").Value = [double]"0.002799731840346333"
$ws.Range("J6").Value = [double]"0.002799731840346333"
$ws.Range("O6").Value = [double]"0.998256289001958"
$ws.Range("P6").Value = [double]"0.998256289001958"
$ws.Range("Q6").Value = [double]"36.46110047306266"
$ws.Range("R6").Value = [double]"328.149904257564"
$ws.Range("S6").Value = [double]"0.002794849917144753"
$ws.Range("T6").Value = [double]"0.002794849917144753"
$ws.Range("I7").Value = [double]"0.002799731840346333"
$ws.Range("J7").Value = [double]"0.002799731840346333"
$ws.Range("M7").Value = [double]"0.2067996666666667"
$ws.Range("N7").Value = [double]"0.620399"
$ws.Range("O7").Value = [double]"0.000762340727132399"
$ws.Range("P7").Value = [double]"0.0007623407271323989"
$ws.Range("Q7").Value = [double]"0.02784433431866667"
$ws.Range("S7").Value = [double]"2.134349606945353E-06"
$ws.Range("T7").Value = [double]"2.134349606945353E-06"
$ws.Range("I8").Value = [double]"0.0382392832655437"
$ws.Range("J8").Value = [double]"0.0382392832655437"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"0.2662156666666667"
$ws.Range("N8").Value = [double]"0.7986470000000001"
$ws.Range("O8").Value = [double]"0.0009813702709097034"
$ws.Range("P8").Value = [double]"0.0009813702709097034"
$ws.Range("Q8").Value = [double]"0.4895690137060001"
$ws.Range("R8").Value = [double]"4.406121123354001"
$ws.Range("S8").Value = [double]"3.752689577769951E-05"
$ws.Range("T8").Value = [double]"3.752689577769951E-05"
$ws.Range("I9").Value = [double]"0.0382392832655437"
$ws.Range("J9").Value = [double]"0.0382392832655437"
$ws.Range("O9").Value = [double]"0.998256289001958"
$ws.Range("P9").Value = [double]"0.998256289001958"
$ws.Range("S9").Value = [double]"0.03817260500675632"
$ws.Range("T9").Value = [double]"0.03817260500675632"
$ws.Range("I10").Value = [double]"0.0382392832655437"
$ws.Range("J10").Value = [double]"0.0382392832655437"
$ws.Range("M10").Value = [double]"0.2067996666666667"
$ws.Range("N10").Value = [double]"0.620399"
$ws.Range("O10").Value = [double]"0.000762340727132399"
$ws.Range("P10").Value = [double]"0.0007623407271323989"
$ws.Range("Q10").Value = [double]"0.3803033462020001"
$ws.Range("R10").Value = [double]"3.422730115818001"
$ws.Range("S10").Value = [double]"2.915136300967636E-05"
$ws.Range("T10").Value = [double]"2.915136300967635E-05"